# Regenerate s_val data to filter save games.
# Updates columns B, C, D, E, G for rows 2-15 on the active sheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$newValues = @{
    2  = @{ B = 1.459612070389937;  C = 1.667794583268128;  D = 0.8054896365839992;  E = 8.660232485948974;  G = 12.59312877619104 }
    3  = @{ B = 0.3048080303191223; C = 0.3127903958511391; D = 3.900430680208489;   E = 8.660232485948974;  G = 13.17826159232772 }
    4  = @{ B = 1.459612070389937;  C = 1.667794583268128;  D = 0.8054896365839992;  E = 8.660232485948974;  G = 12.59312877619104 }
    5  = @{ B = 3.230985683306322;  C = 1.667794583268128;  D = 0.1575252929769615;  E = 0.496779210170732;  G = 5.553084769722144 }
    6  = @{ B = 3.230985683306322;  C = 1.667794583268128;  D = 0.8054896365839992;  E = 0.496779210170732;  G = 6.201049113329182 }
    7  = @{ B = 3.230985683306322;  C = 1.667794583268128;  D = 3.900430680208489;   E = 0.496779210170732;  G = 9.295990156953671 }
    8  = @{ B = 3.230985683306322;  C = 1.667794583268128;  D = 3.900430680208489;   E = 0.496779210170732;  G = 9.295990156953671 }
    9  = @{ B = 3.230985683306322;  C = 1.667794583268128;  D = 3.900430680208489;   E = 0.496779210170732;  G = 9.295990156953671 }
    10 = @{ B = 0.6753301551942219; C = 1.667794583268128;  D = 0.8054896365839992;  E = 0.496779210170732;  G = 3.645393585217082 }
    11 = @{ B = 0.6753301551942219; C = 1.667794583268128;  D = 0.8054896365839992;  E = 0.496779210170732;  G = 3.645393585217082 }
    12 = @{ B = 3.230985683306322;  C = 1.667794583268128;  D = 0.1575252929769615;  E = 0.496779210170732;  G = 5.553084769722144 }
    13 = @{ B = 3.230985683306322;  C = 1.667794583268128;  D = 26.21740644021617;   E = 0.496779210170732;  G = 31.61296591696135 }
    14 = @{ B = 3.230985683306322;  C = 1.667794583268128;  D = 0.8054896365839992;  E = 0.496779210170732;  G = 6.201049113329182 }
    15 = @{ B = 1.459612070389937;  C = 1.667794583268128;  D = 0.8054896365839992;  E = 0.496779210170732;  G = 4.429675500412797 }
}

foreach ($row in $newValues.Keys) {
    $cols = $newValues[$row]
    $ws.Range("B$row").Value = $cols.B
    $ws.Range("C$row").Value = $cols.C
    $ws.Range("D$row").Value = $cols.D
    $ws.Range("E$row").Value = $cols.E
    $ws.Range("G$row").Value = $cols.G
}
